# Update album price list (column D) with refreshed market prices.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 10.33
$ws.Range("D3").Value = 8.779999999999999
$ws.Range("D4").Value = 12.39
$ws.Range("D5").Value = 15.87
$ws.Range("D7").Value = 8.460000000000001
$ws.Range("D9").Value = 5.92
$ws.Range("D10").Value = 24.02
$ws.Range("D11").Value = 1.26
$ws.Range("D12").Value = 2.84
$ws.Range("D13").Value = 1.24
$ws.Range("D14").Value = 6.4
$ws.Range("D15").Value = 3.75
$ws.Range("D16").Value = 4.4
$ws.Range("D17").Value = 16.54
$ws.Range("D18").Value = 3.5
$ws.Range("D19").Value = 12.75
$ws.Range("D20").Value = 4.83
$ws.Range("D21").Value = 1.47
$ws.Range("D22").Value = 6.84
$ws.Range("D23").Value = 2.91
$ws.Range("D24").Value = 3.59
$ws.Range("D26").Value = 4.95
$ws.Range("D28").Value = 4.37
$ws.Range("D30").Value = 0.96
$ws.Range("D31").Value = 2.08
$ws.Range("D32").Value = 1.43
$ws.Range("D33").Value = 3.2
$ws.Range("D36").Value = 1.79
$ws.Range("D37").Value = 5.38
$ws.Range("D38").Value = 25.98
$ws.Range("D39").Value = 0.93
$ws.Range("D40").Value = 1.11
$ws.Range("D42").Value = 3.63
$ws.Range("D43").Value = 3.18
$ws.Range("D44").Value = 0.66
$ws.Range("D46").Value = 3.44
$ws.Range("D47").Value = 8.380000000000001
$ws.Range("D48").Value = 4.35
$ws.Range("D49").Value = 1.22
$ws.Range("D50").Value = 1.78
$ws.Range("D51").Value = 6.07
$ws.Range("D52").Value = 2.02
$ws.Range("D53").Value = 3.54
$ws.Range("D54").Value = 3.57
$ws.Range("D55").Value = 1.44
$ws.Range("D56").Value = 10.99
$ws.Range("D57").Value = 8.58
$ws.Range("D58").Value = 16.33
$ws.Range("D59").Value = 1.34
$ws.Range("D60").Value = 8.859999999999999
$ws.Range("D61").Value = 9.029999999999999
$ws.Range("D62").Value = 3.52
$ws.Range("D63").Value = 1.83
$ws.Range("D64").Value = 5.6
$ws.Range("D65").Value = 3.61
$ws.Range("D66").Value = 1.01
